# Update cryptocurrency price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.796.61'
$ws.Range('E2').Value = '  +0.58%  '

$ws.Range('D3').Value = '3.413.43'
$ws.Range('E3').Value = '  +2.62%  '

$ws.Range('E4').Value = '  -0.27%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '258.86'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.87%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '669.26'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +7.48%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.54'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +7.00%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.469'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +16.94%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.08'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +19.60%  '

$ws.Range('E10').Value = '  -0.17%  '

$ws.Range('D11').Value = '3.407.85'
$ws.Range('E11').Value = '  +2.64%  '

$ws.Range('E12').Value = '  +10.12%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '42.87'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +12.73%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000277'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +10.68%  '

$ws.Range('D15').Value = '98.583.99'
$ws.Range('E15').Value = '  +0.64%  '

$ws.Range('B16').Value = 'Toncoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.81'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +5.31%  '

$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = '4.048.95'
$ws.Range('E17').Value = '  +2.12%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.00'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +31.26%  '

$ws.Range('D19').Value = '3.411.49'
$ws.Range('E19').Value = '  +2.73%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.34'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +14.03%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '532.55'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +10.66%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.56'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.12%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.61'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +12.83%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000219'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +7.01%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.436'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +49.27%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.37'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +13.99%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '102.62'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +15.17%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.78'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +7.44%  '

$ws.Range('D29').Value = '3.596.64'
$ws.Range('E29').Value = '  +2.40%  '

$ws.Range('E30').Value = '  +14.34%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.54'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +17.54%  '

$ws.Range('E32').Value = '  -0.31%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.196'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.98%  '

$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.58%  '

$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '30.47'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +9.43%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.558'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +21.58%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.19'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +12.86%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.163'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +10.60%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.83'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +9.17%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '528.89'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.31%  '

$ws.Range('E41').Value = '  +8.63%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '24.77'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.45%  '

$ws.Range('E43').Value = '  +4.33%  '

$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.55'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +8.47%  '

$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0431'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +32.50%  '

$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.853'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +7.36%  '

$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.12'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +10.30%  '

$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.11'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +16.51%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.21'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +11.93%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.56'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +14.92%  '
